$wb = $excel.ActiveWorkbook

# Update status text: "Ready for handoff" -> "In Translation" on every sheet
# (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4 all carry this status string)
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation") | Out-Null
}

# Narrow the "handoff datetime" columns (report regenerated with a tighter layout)
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5   # column E
$overview.Columns.Item(6).ColumnWidth = 12.5   # column F

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5        # column C

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5        # column C
